# Applies the cryptocurrency price / 1h-volume refresh captured by the
# GitHub Actions scraper run on Sat Apr 29 10:15:42 UTC 2023.
#
# All of the Coin/Link/Price/Volume columns on the sheet are stored as
# plain text (inline strings) in the source workbook, so every touched
# cell is switched to Text number-format ("@") before its value is
# written. That stops Excel from silently re-interpreting values such
# as "1.015" or "5.630" as numbers (which would lose the trailing zero
# / change the stored type) when assigned through .Value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "29.488.08" },
    @{ Cell = "E2"; Value = "  +0.81%  " },
    @{ Cell = "D3"; Value = "1.913.22" },
    @{ Cell = "E3"; Value = "  +0.14%  " },
    @{ Cell = "E4"; Value = "  +0.54%  " },
    @{ Cell = "E5"; Value = "  +1.45%  " },
    @{ Cell = "E6"; Value = "  +0.54%  " },
    @{ Cell = "D7"; Value = "0.4822" },
    @{ Cell = "E7"; Value = "  +2.15%  " },
    @{ Cell = "D8"; Value = "0.4066" },
    @{ Cell = "E8"; Value = "  +0.02%  " },
    @{ Cell = "D9"; Value = "0.08166" },
    @{ Cell = "E9"; Value = "  +1.55%  " },
    @{ Cell = "D10"; Value = "1.015" },
    @{ Cell = "E10"; Value = "  +1.48%  " },
    @{ Cell = "D11"; Value = "23.44" },
    @{ Cell = "E11"; Value = "  +4.30%  " },
    @{ Cell = "D12"; Value = "1.920.27" },
    @{ Cell = "E12"; Value = "  -0.18%  " },
    @{ Cell = "D13"; Value = "6.019" },
    @{ Cell = "E13"; Value = "  +2.36%  " },
    @{ Cell = "D14"; Value = "7.176" },
    @{ Cell = "E14"; Value = "  +0.80%  " },
    @{ Cell = "D15"; Value = "90.34" },
    @{ Cell = "E15"; Value = "  +0.82%  " },
    @{ Cell = "D16"; Value = "0.06798" },
    @{ Cell = "E16"; Value = "  +2.50%  " },
    @{ Cell = "E17"; Value = "  +0.62%  " },
    @{ Cell = "E18"; Value = "  +0.90%  " },
    @{ Cell = "D19"; Value = "17.72" },
    @{ Cell = "E19"; Value = "  +0.38%  " },
    @{ Cell = "E20"; Value = "  +0.51%  " },
    @{ Cell = "D21"; Value = "29.511.59" },
    @{ Cell = "E21"; Value = "  +0.86%  " },
    @{ Cell = "D22"; Value = "5.630" },
    @{ Cell = "E22"; Value = "  +2.08%  " },
    @{ Cell = "E23"; Value = "  +2.64%  " },
    @{ Cell = "D24"; Value = "2.183" },
    @{ Cell = "E24"; Value = "  -0.69%  " },
    @{ Cell = "D25"; Value = "2.151.34" },
    @{ Cell = "E25"; Value = "  +0.15%  " },
    @{ Cell = "D26"; Value = "155.83" },
    @{ Cell = "E26"; Value = "  +0.39%  " },
    @{ Cell = "D27"; Value = "6.447" },
    @{ Cell = "E27"; Value = "  +7.37%  " },
    @{ Cell = "D28"; Value = "20.06" },
    @{ Cell = "E28"; Value = "  +1.49%  " },
    @{ Cell = "D29"; Value = "2.114" },
    @{ Cell = "E29"; Value = "  +0.51%  " },
    @{ Cell = "D30"; Value = "120.22" },
    @{ Cell = "E30"; Value = "  +2.42%  " },
    @{ Cell = "D31"; Value = "1.026" },
    @{ Cell = "E31"; Value = "  -3.68%  " },
    @{ Cell = "D32"; Value = "0.09541" },
    @{ Cell = "E32"; Value = "  +0.40%  " },
    @{ Cell = "D33"; Value = "5.529" },
    @{ Cell = "E33"; Value = "  +2.81%  " },
    @{ Cell = "D34"; Value = "3.566" },
    @{ Cell = "E34"; Value = "  +0.62%  " },
    @{ Cell = "D35"; Value = "1.393" },
    @{ Cell = "E35"; Value = "  -1.95%  " },
    @{ Cell = "D36"; Value = "0.02273" },
    @{ Cell = "E36"; Value = "  +1.47%  " },
    @{ Cell = "D37"; Value = "0.06103" },
    @{ Cell = "E37"; Value = "  +0.50%  " },
    @{ Cell = "D38"; Value = "1.177" },
    @{ Cell = "E38"; Value = "  +0.50%  " },
    @{ Cell = "B39"; Value = "Aptos" },
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt" },
    @{ Cell = "D39"; Value = "10.85" },
    @{ Cell = "E39"; Value = "  +7.39%  " },
    @{ Cell = "B40"; Value = "TheSandbox" },
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand" },
    @{ Cell = "D40"; Value = "0.5963" },
    @{ Cell = "E40"; Value = "  +1.97%  " },
    @{ Cell = "D41"; Value = "8.014" },
    @{ Cell = "E41"; Value = "  -2.47%  " },
    @{ Cell = "D42"; Value = "0.1856" },
    @{ Cell = "E42"; Value = "  +1.26%  " },
    @{ Cell = "E43"; Value = "  +0.61%  " },
    @{ Cell = "D44"; Value = "2.387" },
    @{ Cell = "E44"; Value = "  -4.55%  " },
    @{ Cell = "D45"; Value = "12.54" },
    @{ Cell = "E45"; Value = "  +3.59%  " },
    @{ Cell = "D46"; Value = "0.07611" },
    @{ Cell = "E46"; Value = "  -3.23%  " },
    @{ Cell = "D47"; Value = "0.5579" },
    @{ Cell = "E47"; Value = "  +1.13%  " },
    @{ Cell = "D48"; Value = "1.945" },
    @{ Cell = "E48"; Value = "  +1.42%  " },
    @{ Cell = "D49"; Value = "116.25" },
    @{ Cell = "E49"; Value = "  +2.87%  " },
    @{ Cell = "D50"; Value = "72.68" },
    @{ Cell = "E50"; Value = "  +1.85%  " },
    @{ Cell = "D51"; Value = "2.409" },
    @{ Cell = "E51"; Value = "  +2.72%  " }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    $r.NumberFormat = "@"
    $r.Value = $u.Value
}
